$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 02:12"

# --- Country-name column (A) updates: shared-string reorder shifted which
#     name lands on which already-sorted row; reproduce the end result by
#     writing the new label directly onto the affected rows. ---
$countryNames = @{
    80 = "Bulgaria"
    81 = "Estado de Palestina"
    82 = "Republica de Macedonia"
    102 = "Paraguay"
    103 = "Grecia"
    115 = "Montenegro"
    116 = "Mali"
    138 = "Surinam"
    139 = "Letonia"
    140 = "Uruguay"
    141 = "Jordania"
    142 = "Niger"
    143 = "Liberia"
    171 = "Bahamas"
    172 = "Martinica"
    173 = "Eritrea"
    210 = "Groenlandia"
    211 = "Islas Malvinas"
}
foreach ($r in $countryNames.Keys) {
    $ws.Cells.Item($r, 1).Value = $countryNames[$r]
}

# --- Updated statistics (columns B:H) for the affected rows ---
# Column order per row hashtable: B, C, D, E, F, G, H (use $null to skip/leave unchanged)
$stats = @{
    4 = @(4169146, 68271, 1978873, 2042978, $null, 1112, 147295)
    5 = @(2289951, 58080, 1570237, 635507, $null, 1317, 84207)
    65 = @(18868, 489, 10149, 8617, $null, $null, $null)
    80 = @(9853, 269, 5031, 4493, $null, 8, 329)
    81 = @(9744, 346, 2720, 6957, $null, 1, 67)
    82 = @(9669, 122, 5071, 4153, $null, 3, 445)
    102 = @(4113, 113, 2487, 1590, $null, 0, 36)
    103 = @(4110, 33, 1374, 2535, $null, 1, 201)
    104 = @(3789, 206, $null, 1978, $null, 6, 134)
    109 = @(3171, 10, 1499, 1579, $null, $null, $null)
    115 = @(2569, 97, 538, 1991, $null, 1, 40)
    116 = @(2494, $null, 1889, 482, $null, $null, 123)
    138 = @(1234, 58, 774, 437, $null, 2, 23)
    139 = @(1203, 6, 1045, 127, $null, $null, 31)
    140 = @(1141, 24, 940, 167, $null, $null, 34)
    141 = @(1131, 11, 1035, 85, $null, $null, 11)
    142 = @(1124, 2, 1022, 33, $null, 0, 69)
    143 = @($null, 3, 613, 433, $null, 1, 71)
    164 = @(351, 1, 176, 156, $null, $null, $null)
    171 = @(274, 55, 91, 172, $null, $null, 11)
    172 = @(262, 0, 98, 149, $null, $null, 15)
    173 = @(261, 10, 189, 72, $null, $null, 0)
    181 = @($null, $null, 128, 5, $null, $null, $null)
}
foreach ($r in $stats.Keys) {
    $vals = $stats[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        if ($null -ne $vals[$i]) {
            $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
        }
    }
}
